$d = $word.ActiveDocument

function Add-ProfPrefix($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1)  # exclude the trailing paragraph mark
    $start = $r.Start
    $ins = $d.Range($start, $start)
    $ins.InsertBefore("Prof. ")
    # Force the inserted text into its own run (matching the original run's
    # character formatting exactly) by toggling a character property off/on.
    $profRange = $d.Range($start, $start + 6)
    $profRange.Bold = 1
    $profRange.Bold = 0
}

# "Vishwesha Guttal" -> "Prof. Vishwesha Guttal" (two runs)
Add-ProfPrefix 12
# "Rohini Balakrishnan" -> "Prof. Rohini Balakrishnan" (two runs)
Add-ProfPrefix 13
# "Sutirth Dey" -> "Prof. Sutirth Dey" (two runs)
Add-ProfPrefix 14

# Give the "Sutirth Dey" paragraph explicit before/after spacing (matches the
# spacing that used to live on the trailing empty paragraph being removed).
$sutirth = $d.Paragraphs.Item(14)
$sutirth.Format.SpaceBefore = 0
$sutirth.Format.SpaceAfter = 8

# Remove the now-redundant trailing empty paragraph entirely.
$trailing = $d.Paragraphs.Item(15)
$trailing.Range.Delete()
